# Updates to finalize model
# Adds a third results block ("Starting Algorithm") below the existing
# "Nelder-Mead" and "Quantum Basin Hopping" blocks, and tidies up a few
# leftover/inconsistent cell formats along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------
# 1) Tidy up pre-existing cell formatting that had drifted from the rest
#    of the sheet (center/center alignment everywhere, no stray bold /
#    stray number-format flags).
# ---------------------------------------------------------------------

# Row 2 helper cells: make them fully center/center like the rest of the sheet
$ws.Range("D2:G2").HorizontalAlignment = $xlCenter
$ws.Range("D2:G2").VerticalAlignment = $xlCenter

# "Time" row labels (B11, B24) lost their stray bold/font override
$ws.Range("B11").Font.Bold = $false
$ws.Range("B11").HorizontalAlignment = $xlCenter
$ws.Range("B11").VerticalAlignment = $xlCenter

$ws.Range("B24").Font.Bold = $false
$ws.Range("B24").HorizontalAlignment = $xlCenter
$ws.Range("B24").VerticalAlignment = $xlCenter

# "Time" row values (C11:E11, C24:E24) lost their stray number-format override
$ws.Range("C11:E11").HorizontalAlignment = $xlCenter
$ws.Range("C11:E11").VerticalAlignment = $xlCenter

$ws.Range("C24:E24").HorizontalAlignment = $xlCenter
$ws.Range("C24:E24").VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------
# 2) Add the new "Starting Algorithm" results block in rows 29-39,
#    mirroring the layout of the two blocks above it.
# ---------------------------------------------------------------------

# Merged section header
$ws.Range("C29:E29").Merge()
$ws.Range("C29").Value = "Starting Algorithm"
$ws.Range("C29:E29").HorizontalAlignment = $xlCenter
$ws.Range("C29:E29").VerticalAlignment = $xlCenter

# Trials / Shots / Iterations / Seeds / Nodes / Edges
$ws.Range("B30").Value = "Trials"
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 3

$ws.Range("B31").Value = "Shots"
$ws.Range("C31").Value = 5000
$ws.Range("D31").Value = 5000
$ws.Range("E31").Value = 5000

$ws.Range("B32").Value = "Iterations"
$ws.Range("C32").Value = 100
$ws.Range("D32").Value = 100
$ws.Range("E32").Value = 100

$ws.Range("B33").Value = "Seeds"
$ws.Range("C33").Value = 12345
$ws.Range("D33").Value = 12345
$ws.Range("E33").Value = 12345

$ws.Range("B34").Value = "Nodes"
$ws.Range("C34").Value = 7
$ws.Range("D34").Value = 7
$ws.Range("E34").Value = 7

$ws.Range("B35").Value = "Edges"
$ws.Range("C35").Value = 6
$ws.Range("D35").Value = 6
$ws.Range("E35").Value = 6

$ws.Range("B30:E35").HorizontalAlignment = $xlCenter
$ws.Range("B30:E35").VerticalAlignment = $xlCenter

# "Outputs" sub-header (bold label, like rows 10 and 23)
$ws.Range("B36").Value = "Outputs"
$ws.Range("B36").Font.Bold = $true
$ws.Range("B36").HorizontalAlignment = $xlCenter
$ws.Range("B36").VerticalAlignment = $xlCenter

# Time / Higest Energy / Success Rate (single trial, only column C populated)
$ws.Range("B37").Value = "Time"
$ws.Range("C37").Value = 95

$ws.Range("B38").Value = "Higest Energy"
$ws.Range("C38").Value = 4.94

$ws.Range("B39").Value = "Success Rate"
$ws.Range("C39").Value = 0.42
$ws.Range("C39").NumberFormat = "0%"

$ws.Range("B37:C39").HorizontalAlignment = $xlCenter
$ws.Range("B37:C39").VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------
# 3) Update the active selection / scroll position to match where the
#    edit was made.
# ---------------------------------------------------------------------
$win = $excel.Windows.Item(1)
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("D37").Select()
